{"js": "// 1) Header-block paragraph in the document body:\n//    \"Title: Prototyping Labs Manager\" -> \"Title: Prototyping Lab Manager\"\n{\n  const results = context.document.body.search(\"Prototyping Labs Manager\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"Prototyping Lab Manager\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 2) Primary (default) page header: \"Prototyping Labs at GIX\" -> \"Prototyping Lab at GIX\"\n{\n  const sections = context.document.sections;\n  sections.load(\"items\");\n  await context.sync();\n\n  for (let s = 0; s < sections.items.length; s++) {\n    const header = sections.items[s].getHeader(Word.HeaderFooterType.primary);\n    const results = header.search(\"Prototyping Labs at GIX\", { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n      results.items[i].insertText(\"Prototyping Lab at GIX\", Word.InsertLocation.replace);\n    }\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Header-block paragraph in the document body:\n#    \"Title: Prototyping Labs Manager\" -> \"Title: Prototyping Lab Manager\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Prototyping Labs Manager\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Prototyping Lab Manager\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\n# 2) Primary (default) page header: \"Prototyping Labs at GIX\" -> \"Prototyping Lab at GIX\"\nforeach ($sec in $d.Sections) {\n    $hdr = $sec.Headers.Item(1)  # wdHeaderFooterPrimary\n    $hfind = $hdr.Range.Find\n    $hfind.ClearFormatting()\n    $hfind.Text = \"Prototyping Labs at GIX\"\n    $hfind.Replacement.ClearFormatting()\n    $hfind.Replacement.Text = \"Prototyping Lab at GIX\"\n    $hfind.Execute([ref]$hfind.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$hfind.Replacement.Text, 2)\n}\n"}
